# Delete the row for TSN (Tianjin, China) from the colo data table.
# This shifts all subsequent rows up by one, removing the final row
# (YHZ / Halifax, Canada) and shrinking the used range from A1:H334
# to A1:H333.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(235).Delete()
